# Apply the OOXML diff to the workbook using Excel COM-interop semantics.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the shared-string text for the "probability / z-score" question.
#    In the "before" workbook this text lives on sheet "7_" (A1).  We update
#    the cell value in place; Excel will drop the now-unused string and
#    append a new one, which is exactly what the target diff shows
#    (old si 98 removed, new si 112 appended at the very end).
# ---------------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item("7_")
$ws7.Range("A1").Value = 'You take data on the average voltage of a set of batteries, and find an average value of V=9.15 volts, with a standard deviation of 0.36.  What is the probability that the next battery you check will have a value within 0.5 volts of 9.15?  Use a z-score calculator to find the answer (see the link above).  Enter your answer as a whole number percentage (for example, "12" for 12.3%)'
$ws7.Rows.Item(1).RowHeight = 180

# ---------------------------------------------------------------------------
# 2. On sheet "10_" swap which row holds the "correct" answer: the Y/N
#    markers and the explanation text move from row 3 (0.09) to row 4 (0.18).
# ---------------------------------------------------------------------------
$ws10b = $wb.Worksheets.Item("10_")
$ws10b.Range("B3").Value = "N"
$ws10b.Range("C3").Value = ""
$ws10b.Rows.Item(3).AutoFit()

$ws10b.Range("B4").Value = "Y"
$ws10b.Range("C4").Value = "The standard deviation of the mean is the Sdev divided by the square root of 9.  So SMD = 0.36/4 = 0.9.   A 95% confidence level requires a z-score of 2.0, and uncertainty =z*SDM = 0.18."
$ws10b.Rows.Item(4).RowHeight = 60

# ---------------------------------------------------------------------------
# 3. Selection / active-cell bookkeeping on the various sheets.  These are
#    applied in order, and whichever sheet is activated/selected last ends
#    up as the workbook's active tab, matching the target (sheet "10_").
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("4_")
$ws4.Range("D2").Select()

$ws7.Range("F4").Select()

# "10_" is activated last, so it becomes the workbook's active sheet/tab,
# matching activeTab="10" and tabSelected="1" in the target workbook.
$ws10b.Range("C1").Select()
